$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 11 (dimension A1:R11). A new weekly
# price record is being appended as row 12, carrying the values that
# used to live in row 11, while row 11 itself is updated in place with
# the new week's figures.

# Duplicate row 11 (current/old values) down into row 12 first.
$oldRow = $ws.Rows.Item(11).Value2
$ws.Rows.Item(12).Value2 = $oldRow

# Row 12's date cell (D12) should keep the same date number format as D11.
$ws.Range("D11").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Now overwrite row 11 with the new week's values per the update.
$ws.Range("D11").Value = 44476
$ws.Range("J11").Value = 160
$ws.Range("K11").Value = 7500
$ws.Range("L11").Value = 8000
$ws.Range("M11").Value = 7750
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 310
